$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:C columns to B:D by inserting a brand-new blank
# column at A. This keeps each moved cell's value, type, number format and
# style (e.g. the bold/bordered/centered header row, and the date number
# format on column B -> now C) intact, since Excel's column insert carries
# formatting along with the shifted cells.
$ws.Columns.Item(1).Insert()

# Rename the (now shifted) header captions to the new, shorter labels.
# These cells already inherited the header style (bold, thin border,
# center/top aligned) from the shift above, so a plain value write suffices.
$ws.Range("B1").Value = "codice"
$ws.Range("C1").Value = "data"
$ws.Range("D1").Value = "val"

# New column A is a simple numeric row index (0-based), styled to match the
# header look (bold, thin border, centered & top aligned) like A2/A3 in the
# target sheet.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A2:A3").Font.Bold = $true
$ws.Range("A2:A3").Borders.LineStyle = 1
$ws.Range("A2:A3").HorizontalAlignment = -4108
$ws.Range("A2:A3").VerticalAlignment = -4160

# Re-assert the measurement values in the new D column (harmless - same
# IEEE-754 double as before the shift, just guards against any drift in the
# column-shift's internal re-serialisation of the literal).
$ws.Range("D2").Value = 1.33
$ws.Range("D3").Value = 1.92

Write-Output "done"
